$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-04 Saturday" "2024-05-05 Sunday"

Replace-Text "62×61=3782" "93×75=6975"
Replace-Text "54×68=3672" "12×39=468"
Replace-Text "64×35=2240" "72×68=4896"
Replace-Text "70×56=3920" "87×15=1305"
Replace-Text "29×30=870" "20×96=1920"
Replace-Text "98×66=6468" "54×12=648"
Replace-Text "45×41=1845" "78×96=7488"
Replace-Text "41×98=4018" "21×97=2037"
Replace-Text "22×93=2046" "75×64=4800"
Replace-Text "81×79=6399" "52×86=4472"
Replace-Text "71×63=4473" "43×87=3741"
Replace-Text "99×19=1881" "44×67=2948"
Replace-Text "25×17=425" "21×96=2016"
Replace-Text "83×71=5893" "54×41=2214"
Replace-Text "44×85=3740" "62×98=6076"
Replace-Text "25×34=850" "15×47=705"
Replace-Text "90×96=8640" "81×95=7695"
Replace-Text "47×18=846" "30×27=810"
Replace-Text "43×60=2580" "66×66=4356"
Replace-Text "77×60=4620" "53×83=4399"
Replace-Text "53×12=636" "13×82=1066"
Replace-Text "74×93=6882" "66×99=6534"
Replace-Text "98×82=8036" "88×66=5808"
Replace-Text "72×98=7056" "99×17=1683"
Replace-Text "80×98=7840" "42×27=1134"
